$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$doi = "DOI: 10.1016/j.jchromb.2018.01.035 "
$doiUrl = "https://doi.org/10.1016/j.jchromb.2018.01.035"

# --- Row 4: Heptane ---
$ws.Range("A4").Value = "Heptane"
$ws.Range("B4").Value = 15.3
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = "HSP Handbook"

# --- Row 5: H2O ---
$ws.Range("A5").Value = "H2O"
$ws.Range("B5").Value = 15.5
$ws.Range("C5").Value = 16
$ws.Range("D5").Value = 42.3
$ws.Range("E5").Value = "HSP Handbook"

# --- Row 6: Luteolin (own hyperlink) ---
[void]$ws.Hyperlinks.Add($ws.Range("E6"), $doiUrl, "", "", $doiUrl)
$ws.Range("E6").Value = $doi
$ws.Range("A6").Value = "Luteolin"
$ws.Range("B6").Value = 20.6
$ws.Range("C6").Value = 9.1
$ws.Range("D6").Value = 10.8

# --- Rows 7-9: Wedelolactone, Apigenin, Quercetin (shared hyperlink E7:E9) ---
[void]$ws.Hyperlinks.Add($ws.Range("E7:E9"), $doiUrl, "", "", $doiUrl)
$ws.Range("E7:E9").Value = $doi
$ws.Range("E7:E9").Style = "Hyperlink"

$ws.Range("A7").Value = "Wedelolactone"
$ws.Range("B7").Value = 20.4
$ws.Range("C7").Value = 7.5
$ws.Range("D7").Value = 10.8

$ws.Range("A8").Value = "Apigenin"
$ws.Range("B8").Value = 20.4
$ws.Range("C8").Value = 9.2
$ws.Range("D8").Value = 13.5

$ws.Range("A9").Value = "Quercetin"
$ws.Range("B9").Value = 21
$ws.Range("C9").Value = 10.6
$ws.Range("D9").Value = 13.7

# --- Column E width ---
$ws.Range("E1").ColumnWidth = 30.75

# --- Selection ---
[void]$ws.Range("E6").Select()
